# Apply "Doing Updates for Financials" changes to the GVP worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVP")

# Non Recurring (row 21) - FY2012 (col J) value became unavailable -> "NA"
$ws.Range("J21").Value = "NA"

# Depreciation (row 83) - FY2012 (col J) value became unavailable -> "NA"
$ws.Range("J83").Value = "NA"

# Other Cashflows from Investing Activities (row 91) - values restated
$ws.Range("D91").Value = -100
$ws.Range("E91").Value = -100
$ws.Range("F91").Value = -300
$ws.Range("G91").Value = -400
$ws.Range("H91").Value = -400
$ws.Range("I91").Value = -1600
$ws.Range("J91").Value = -500

# Total Cash Flows From Investing Activities (row 94) - FY2012 (col J) value became unavailable -> "NA"
$ws.Range("J94").Value = "NA"

# Other Cash Flows from Financing Activities (row 100) - FY2012 (col J) value became unavailable -> "NA"
$ws.Range("J100").Value = "NA"

# Total Cash Flows From Financing Activities (row 101) - FY2012 (col J) value became unavailable -> "NA"
$ws.Range("J101").Value = "NA"
